$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Change 1: delete the standalone paragraph "Research question(s)."
# that duplicates the "Section 1 - Research Question(s)" heading.
# ---------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Research question(s).`r") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------
# Change 2: join the paragraph ending in "...receive." with the
# following paragraph "This will be attempted using..." (delete the
# paragraph mark between them, no extra space introduced).
# ---------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*level of Sisben students*parents receive.`r") {
        $mark = $d.Range($p.Range.End - 1, $p.Range.End)
        $mark.Delete()
        break
    }
}

# ---------------------------------------------------------------
# Change 3: "... are all binary " -> "... are all binary." and add a
# new paragraph after it: "Sisben is a discrete ordinal variable."
# (matching the surrounding minorHAnsi / 2E2E2E formatting).
# ---------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = " are all binary "
$find.Replacement.ClearFormatting()
$find.Replacement.Text = " are all binary."
$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*MOBILE are all binary.`r") {
        # Nudge the trailing "binary." into its own run (mirrors the
        # author typing the period right after "binary" as a separate
        # edit) while keeping the inherited rFonts/color formatting.
        $tailLen = "binary.".Length
        $tailRange = $d.Range($p.Range.End - 1 - $tailLen, $p.Range.End - 1)
        $tailRange.Font.Bold = $true
        $tailRange.Font.Bold = $false

        $p.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Range.Text = "Sisben is a discrete ordinal variable."
        break
    }
}

# ---------------------------------------------------------------
# Change 4: "... ENG_S11 are all discrete variables" -> add a
# trailing period: "... ENG_S11 are all discrete variables."
# ---------------------------------------------------------------
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Text = "ENG_S11 are all discrete variables"
$find2.Replacement.ClearFormatting()
$find2.Replacement.Text = "ENG_S11 are all discrete variables."
$find2.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Saber 11 exams:*discrete variables.`r") {
        $dotRange = $d.Range($p.Range.End - 2, $p.Range.End - 1)
        $dotRange.Font.Bold = $true
        $dotRange.Font.Bold = $false
        break
    }
}
